$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sending cluster names (new TPM data reshuffled cluster labels)
$ws.Range("A2").Value = "MuSCs"
$ws.Range("A3").Value = "Resolving-Mac"

# Row 2 - updated ligand expression values (new TPM) and dependent specificity/weight columns
$ws.Range("G2").Value = 0.2468816666666667
$ws.Range("H2").Value = 0.740645
$ws.Range("I2").Value = 0.6299860588115711
$ws.Range("J2").Value = 0.6299860588115711
$ws.Range("Q2").Value = 0.1350167855077778
$ws.Range("R2").Value = 1.21515106957
$ws.Range("S2").Value = 0.6299860588115711
$ws.Range("T2").Value = 0.6299860588115711

# Row 3 - updated ligand expression values (new TPM) and dependent specificity/weight columns
$ws.Range("G3").Value = 0.1450026666666667
$ws.Range("H3").Value = 0.435008
$ws.Range("I3").Value = 0.3700139411884289
$ws.Range("J3").Value = 0.3700139411884289
$ws.Range("Q3").Value = 0.07930031503644445
$ws.Range("R3").Value = 0.7137028353279999
$ws.Range("S3").Value = 0.3700139411884289
$ws.Range("T3").Value = 0.3700139411884289
